$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.815.01'
$ws.Range("E2").Value = '  -0.56%  '
$ws.Range("D3").Value = '2.305.44'
$ws.Range("E3").Value = '  -2.15%  '
$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''544.05'
$ws.Range("E5").Value = '  -0.30%  '
$ws.Range("D6").Value = '''129.17'
$ws.Range("E6").Value = '  -2.35%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  -2.52%  '
$ws.Range("D9").Value = '2.302.84'
$ws.Range("E9").Value = '  -2.21%  '
$ws.Range("E10").Value = '  -0.25%  '
$ws.Range("E11").Value = '  +0.50%  '
$ws.Range("E12").Value = '  -0.65%  '
$ws.Range("E13").Value = '  -0.31%  '
$ws.Range("D14").Value = '''23.38'
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '59.792.53'
$ws.Range("E15").Value = '  -0.42%  '
$ws.Range("B16").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C16").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D16").Value = '2.712.22'
$ws.Range("E16").Value = '  -2.22%  '
$ws.Range("E17").Value = '  +0.01%  '
$ws.Range("D18").Value = '2.310.06'
$ws.Range("E18").Value = '  -1.61%  '
$ws.Range("E19").Value = '  -1.80%  '
$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").Value = '''4.06'
$ws.Range("E20").Value = '  -2.58%  '
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").Value = '''312.00'
$ws.Range("E21").Value = '  -0.68%  '
$ws.Range("E22").Value = '  -4.19%  '
$ws.Range("E23").Value = '  +0.16%  '
$ws.Range("D24").Value = '''64.01'
$ws.Range("E24").Value = '  +1.12%  '
$ws.Range("D25").Value = '''0.170'
$ws.Range("E25").Value = '  -1.21%  '
$ws.Range("D26").Value = '''0.999'
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").Value = '''7.79'
$ws.Range("E27").Value = '  -1.45%  '
$ws.Range("E28").Value = '  +0.23%  '
$ws.Range("D29").Value = '''1.24'
$ws.Range("E29").Value = '  +7.73%  '
$ws.Range("D30").Value = '''170.63'
$ws.Range("E30").Value = '  -0.62%  '
$ws.Range("D31").Value = '''1.71'
$ws.Range("E31").Value = '  -1.64%  '
$ws.Range("D32").Value = '0.0₃0722'
$ws.Range("E32").Value = '  -0.95%  '
$ws.Range("D33").Value = '''5.95'
$ws.Range("E33").Value = '  +0.10%  '
$ws.Range("D34").Value = '''0.380'
$ws.Range("E34").Value = '  -0.62%  '
$ws.Range("E35").Value = '  -4.30%  '
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("D37").Value = '''17.90'
$ws.Range("E37").Value = '  -0.79%  '
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("D39").Value = '''4.05'
$ws.Range("E39").Value = '  -2.72%  '
$ws.Range("D40").Value = '''314.69'
$ws.Range("E40").Value = '  -1.42%  '
$ws.Range("D41").Value = '''37.92'
$ws.Range("E41").Value = '  -0.47%  '
$ws.Range("E42").Value = '  -1.45%  '
$ws.Range("D43").Value = '''136.84'
$ws.Range("E43").Value = '  -4.05%  '
$ws.Range("E44").Value = '  +0.40%  '
$ws.Range("D45").Value = '''0.0938'
$ws.Range("E45").Value = '  -1.38%  '
$ws.Range("D46").Value = '''18.98'
$ws.Range("E46").Value = '  -2.27%  '
$ws.Range("E47").Value = '  -0.36%  '
$ws.Range("E48").Value = '  -1.18%  '
$ws.Range("E49").Value = '  -0.35%  '
$ws.Range("D50").Value = '0.0₆0216'
$ws.Range("E50").Value = '  +2.77%  '
$ws.Range("D51").Value = '''16.79'
$ws.Range("E51").Value = '  -1.23%  '
